$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header update: I1 "Pot_Profit" -> "Profit" ---
$ws.Range("I1").Value = "Profit"

# --- Row 2: Manchester United (HOME) - values updated, identity unchanged ---
$ws.Range("A2").Value = "HOME"
$ws.Range("B2").Value = "Manchester United"
$ws.Range("C2").Value = 2.53
$ws.Range("D2").Value = 0.537014889210818
$ws.Range("E2").Value = 0.3586476697033694
$ws.Range("F2").Value = "Manchester United vs Newcastle United"
$ws.Range("G2").Value = 0.04688204832723783
$ws.Range("H2").Value = 1.406461449817135
$ws.Range("I2").Value = 2.151886018220216

# --- Row 3: now Aston Villa (AWAY), was Burnley ---
$ws.Range("A3").Value = "AWAY"
$ws.Range("B3").Value = "Aston Villa"
$ws.Range("C3").Value = 4.03
$ws.Range("D3").Value = 0.3927878780681108
$ws.Range("E3").Value = 0.5829351486144867
$ws.Range("F3").Value = "Chelsea vs Aston Villa"
$ws.Range("G3").Value = 0.03078205405225012
$ws.Range("H3").Value = 0.9234616215675037
$ws.Range("I3").Value = 2.798088713349536

# --- Row 4: now Burnley (HOME), was Nottingham Forest ---
$ws.Range("A4").Value = "HOME"
$ws.Range("B4").Value = "Burnley"
$ws.Range("C4").Value = 4.01
$ws.Range("D4").Value = 0.3896918319201798
$ws.Range("E4").Value = 0.5626642459999209
$ws.Range("F4").Value = "Burnley vs Everton"
$ws.Range("G4").Value = 0.0299090629102948
$ws.Range("H4").Value = 0.8972718873088441
$ws.Range("I4").Value = 2.70078838079962

# --- Row 5: now West Ham United (HOME), was Aston Villa (AWAY) ---
$ws.Range("A5").Value = "HOME"
$ws.Range("B5").Value = "West Ham United"
$ws.Range("C5").Value = 2.68
$ws.Range("D5").Value = 0.4438756880932579
$ws.Range("E5").Value = 0.1895868440899313
$ws.Range("F5").Value = "West Ham United vs Fulham"
$ws.Range("G5").Value = 0.02256986239165849
$ws.Range("H5").Value = 0.6770958717497547
$ws.Range("I5").Value = 1.137521064539588

# --- Row 6: now Nottingham Forest (HOME), was West Ham United ---
$ws.Range("A6").Value = "HOME"
$ws.Range("B6").Value = "Nottingham Forest"
$ws.Range("C6").Value = 5.29
$ws.Range("D6").Value = 0.2345800676157228
$ws.Range("E6").Value = 0.2409285576871736
$ws.Range("F6").Value = "Nottingham Forest vs Manchester City"
$ws.Range("G6").Value = 0.006739260354885974
$ws.Range("H6").Value = 0.2021778106465792
$ws.Range("I6").Value = 0.8673428076738249

# --- Row 7: now DRAW / Draw (Arsenal vs Brighton), was Brentford (HOME) ---
$ws.Range("A7").Value = "DRAW"
$ws.Range("B7").Value = "Draw"
$ws.Range("C7").Value = 4.92
$ws.Range("D7").Value = 0.2189947313686686
$ws.Range("E7").Value = 0.07745407833384976
$ws.Range("F7").Value = "Arsenal vs Brighton & Hove Albion"
$ws.Range("G7").Value = 0.002371043214301524
$ws.Range("H7").Value = 0.07113129642904571
$ws.Range("I7").Value = 0.2788346820018592
